$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "2025-04-17" row (original row 2); remaining rows shift up.
$ws.Rows.Item(2).Delete()

# Refresh the randomized "quantidade_atipica" index column (A) for the new row set.
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(6, 1).Value = 0
$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(11, 1).Value = 3
